$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the existing (zero-width) "_GoBack" bookmark that currently
#    sits in the empty paragraph right after "Khong co giao dich phat
#    sinh hoac thu chi nen se ko co du lieu".
# ------------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# ------------------------------------------------------------------
# 2. Mark "- Xu ly mail loi: VJ ngay 9/12/20 MBT48R" as done (strike-through).
# ------------------------------------------------------------------
$target1 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Xử lý mail lỗi: VJ ngày 9/12/20 MBT48R*") {
        $target1 = $p
        break
    }
}
$target1.Range.Font.StrikeThrough = 1

# ------------------------------------------------------------------
# 3. Mark "- Ko xu ly dc mail ma chang bay khac SGN-HAN: Mail BB ngay
#    8/12/20" as done too (strike-through) and move the "_GoBack"
#    bookmark so it wraps this paragraph (new last-edit position).
# ------------------------------------------------------------------
$target2 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Ko xử lý đc mail mà chặng bay khác SGN-HAN*") {
        $target2 = $p
        break
    }
}
$target2.Range.Font.StrikeThrough = 1

$d.Bookmarks.Add("_GoBack", $target2.Range)

Write-Output "done"
